# Auto-generated edit script: updates computed market columns (H-N)
# across multiple worksheets, per the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H33").Value = 242.20833
$ws.Range("I33").Value = 222.88889
$ws.Range("J33").Value = 300.16666
$ws.Range("K33").Value = 222.88889
$ws.Range("L33").Value = 300.16666
$ws.Range("M33").Value = 6.111109999999996
$ws.Range("N33").Value = -758.16666

$ws.Range("H62").Value = 2040.4166
$ws.Range("I62").Value = 1840.579
$ws.Range("J62").Value = 2799.8
$ws.Range("K62").Value = 1840.579
$ws.Range("L62").Value = 2799.8
$ws.Range("M62").Value = -1216.579
$ws.Range("N62").Value = -4047.8

$ws.Range("H65").Value = 2040.4166
$ws.Range("I65").Value = 1840.579
$ws.Range("J65").Value = 2799.8
$ws.Range("K65").Value = 9202.895
$ws.Range("L65").Value = 13999
$ws.Range("M65").Value = -6082.895
$ws.Range("N65").Value = -20239

$ws.Range("H121").Value = 705.625
$ws.Range("I121").Value = 595
$ws.Range("J121").Value = 710.43475
$ws.Range("K121").Value = 1785
$ws.Range("L121").Value = 2131.30425
$ws.Range("M121").Value = -38
$ws.Range("N121").Value = -5625.30425

$ws.Range("H129").Value = 880.4815
$ws.Range("I129").Value = 457.44446
$ws.Range("J129").Value = 965.08887
$ws.Range("K129").Value = 1372.33338
$ws.Range("L129").Value = 2895.26661
$ws.Range("M129").Value = 3627.66662
$ws.Range("N129").Value = -12895.26661

$ws.Range("H132").Value = 1636277.1
$ws.Range("I132").Value = 3200.6667
$ws.Range("J132").Value = 16333965
$ws.Range("K132").Value = 9602.000100000001
$ws.Range("L132").Value = 49001895
$ws.Range("M132").Value = -7072.000100000001
$ws.Range("N132").Value = -49006955

$ws.Range("H137").Value = 2780130.5
$ws.Range("I137").Value = 5884105.5
$ws.Range("J137").Value = 2889.4736
$ws.Range("K137").Value = 17652316.5
$ws.Range("L137").Value = 8668.4208
$ws.Range("M137").Value = -17649766.5
$ws.Range("N137").Value = -13768.4208

$ws.Range("H138").Value = 2319778.8
$ws.Range("I138").Value = 3979.4
$ws.Range("J138").Value = 2693294.8
$ws.Range("K138").Value = 11938.2
$ws.Range("L138").Value = 8079884.399999999
$ws.Range("M138").Value = -6798.200000000001
$ws.Range("N138").Value = -8090164.399999999

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H14").Value = 5835
$ws.Range("I14").Value = 5835
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 5835
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -5660

$ws.Range("H32").Value = 8160122
$ws.Range("I32").Value = 9557053
$ws.Range("J32").Value = 11356.833
$ws.Range("K32").Value = 9557053
$ws.Range("L32").Value = 11356.833
$ws.Range("M32").Value = -9556766
$ws.Range("N32").Value = -11930.833

$ws.Range("H74").Value = 5479127
$ws.Range("I74").Value = 9296830
$ws.Range("J74").Value = 53970.21
$ws.Range("K74").Value = 9296830
$ws.Range("L74").Value = 53970.21
$ws.Range("M74").Value = -9295956
$ws.Range("N74").Value = -55718.21

$ws.Range("H77").Value = 5479127
$ws.Range("I77").Value = 9296830
$ws.Range("J77").Value = 53970.21
$ws.Range("K77").Value = 46484150
$ws.Range("L77").Value = 269851.05
$ws.Range("M77").Value = -46479782
$ws.Range("N77").Value = -278587.05

$ws.Range("H88").Value = 5469.44
$ws.Range("I88").Value = 3454
$ws.Range("J88").Value = 6603.125
$ws.Range("K88").Value = 3454
$ws.Range("L88").Value = 6603.125
$ws.Range("M88").Value = -3048
$ws.Range("N88").Value = -7415.125

$ws.Range("H91").Value = 5469.44
$ws.Range("I91").Value = 3454
$ws.Range("J91").Value = 6603.125
$ws.Range("K91").Value = 3454
$ws.Range("L91").Value = 6603.125
$ws.Range("M91").Value = -2050
$ws.Range("N91").Value = -9411.125

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H43").Value = 149980
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 149980
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 149980
$ws.Range("N43").Value = -150342

$ws.Range("H86").Value = 26966.666
$ws.Range("I86").Value = 29975
$ws.Range("J86").Value = 2900
$ws.Range("K86").Value = 29975
$ws.Range("L86").Value = 2900
$ws.Range("M86").Value = -28852
$ws.Range("N86").Value = -5146

$ws.Range("H89").Value = 26966.666
$ws.Range("I89").Value = 29975
$ws.Range("J89").Value = 2900
$ws.Range("K89").Value = 149875
$ws.Range("L89").Value = 14500
$ws.Range("M89").Value = -144259
$ws.Range("N89").Value = -25732

$ws.Range("H134").Value = 2907.3818
$ws.Range("I134").Value = 2727.2195
$ws.Range("J134").Value = 3435
$ws.Range("K134").Value = 8181.6585
$ws.Range("L134").Value = 10305
$ws.Range("M134").Value = -5646.6585
$ws.Range("N134").Value = -15375

$ws.Range("H141").Value = 49000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 49000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 49000
$ws.Range("N141").Value = -59360

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H31").Value = 1228.6842
$ws.Range("I31").Value = 780.0741
$ws.Range("J31").Value = 1406.8088
$ws.Range("K31").Value = 780.0741
$ws.Range("L31").Value = 1406.8088
$ws.Range("M31").Value = -485.0741
$ws.Range("N31").Value = -1996.8088

$ws.Range("H34").Value = 1228.6842
$ws.Range("I34").Value = 780.0741
$ws.Range("J34").Value = 1406.8088
$ws.Range("K34").Value = 780.0741
$ws.Range("L34").Value = 1406.8088
$ws.Range("M34").Value = -578.0741
$ws.Range("N34").Value = -1810.8088

$ws.Range("H132").Value = 62090.53
$ws.Range("I132").Value = 3454.1
$ws.Range("J132").Value = 145856.86
$ws.Range("K132").Value = 10362.3
$ws.Range("L132").Value = 437570.58
$ws.Range("M132").Value = -7832.299999999999
$ws.Range("N132").Value = -442630.58

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H5").Value = 19023.092
$ws.Range("I5").Value = 29769.47
$ws.Range("J5").Value = 754.25
$ws.Range("K5").Value = 89308.41
$ws.Range("L5").Value = 2262.75
$ws.Range("M5").Value = -89196.41
$ws.Range("N5").Value = -2486.75

$ws.Range("H68").Value = 1199.6234
$ws.Range("I68").Value = 651.5769
$ws.Range("J68").Value = 1479.0197
$ws.Range("K68").Value = 1954.7307
$ws.Range("L68").Value = 4437.0591
$ws.Range("M68").Value = -1143.7307
$ws.Range("N68").Value = -6059.0591

$ws.Range("H71").Value = 1199.6234
$ws.Range("I71").Value = 651.5769
$ws.Range("J71").Value = 1479.0197
$ws.Range("K71").Value = 5864.1921
$ws.Range("L71").Value = 13311.1773
$ws.Range("M71").Value = -1808.1921
$ws.Range("N71").Value = -21423.1773

$ws.Range("H107").Value = 1015.3563
$ws.Range("I107").Value = 504.1
$ws.Range("J107").Value = 1450.4681
$ws.Range("K107").Value = 1512.3
$ws.Range("L107").Value = 4351.4043
$ws.Range("M107").Value = 407.6999999999998
$ws.Range("N107").Value = -8191.4043

$ws.Range("H131").Value = 689.4737
$ws.Range("I131").Value = 448
$ws.Range("J131").Value = 957.7778
$ws.Range("K131").Value = 1344
$ws.Range("L131").Value = 2873.3334
$ws.Range("M131").Value = 3696
$ws.Range("N131").Value = -12953.3334

$ws.Range("H135").Value = 19023.092
$ws.Range("I135").Value = 29769.47
$ws.Range("J135").Value = 754.25
$ws.Range("K135").Value = 267925.23
$ws.Range("L135").Value = 6788.25
$ws.Range("M135").Value = -265390.23
$ws.Range("N135").Value = -11858.25

$ws.Range("H139").Value = 3232.75
$ws.Range("I139").Value = 1476.2106
$ws.Range("J139").Value = 5800
$ws.Range("K139").Value = 4428.6318
$ws.Range("L139").Value = 17400
$ws.Range("M139").Value = 711.3681999999999
$ws.Range("N139").Value = -27680

$ws.Range("H140").Value = 1992.68
$ws.Range("I140").Value = 961.3333
$ws.Range("J140").Value = 2174.6824
$ws.Range("K140").Value = 2883.9999
$ws.Range("L140").Value = 6524.047200000001
$ws.Range("M140").Value = 2296.0001
$ws.Range("N140").Value = -16884.0472

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H44").Value = 9964
$ws.Range("I44").Value = 9928
$ws.Range("J44").Value = 10000
$ws.Range("K44").Value = 9928
$ws.Range("L44").Value = 10000
$ws.Range("M44").Value = -9332
$ws.Range("N44").Value = -11192

$ws.Range("H132").Value = 37766.555
$ws.Range("I132").Value = 23878.422
$ws.Range("J132").Value = 94581.63
$ws.Range("K132").Value = 71635.266
$ws.Range("L132").Value = 283744.89
$ws.Range("M132").Value = -69105.266
$ws.Range("N132").Value = -288804.89

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H132").Value = 34886.902
$ws.Range("I132").Value = 1762.0625
$ws.Range("J132").Value = 70220.07000000001
$ws.Range("K132").Value = 5286.1875
$ws.Range("L132").Value = 210660.21
$ws.Range("M132").Value = -2756.1875
$ws.Range("N132").Value = -215720.21

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H113").Value = 1607.125
$ws.Range("I113").Value = 581.9
$ws.Range("J113").Value = 3315.8333
$ws.Range("K113").Value = 1745.7
$ws.Range("L113").Value = 9947.499899999999
$ws.Range("M113").Value = 424.3000000000002
$ws.Range("N113").Value = -14287.4999

$ws.Range("H122").Value = 3912.25
$ws.Range("I122").Value = 3794
$ws.Range("J122").Value = 3929.1428
$ws.Range("K122").Value = 11382
$ws.Range("L122").Value = 11787.4284
$ws.Range("M122").Value = -8932
$ws.Range("N122").Value = -16687.4284
